# Add data for 2022-01-01: update the "through" date in the sheet name and
# in the December row label, and refresh the December + Total rows with the
# new arrest / no-arrest counts (and recomputed arrest_rate percentages).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the worksheet tab (and workbook <sheet name>) ---
$ws.Name = "Through 2021-12-24"

# --- Update the row label for the December row ---
$ws.Range("A14").Value = "December (through 12-24)"

# --- December row (row 14): arrest_made, no_arrest_made, arrest_rate per year ---
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 0.1176

$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 70
$ws.Range("G14").Value = 0.0909

$ws.Range("H14").Value = 10
$ws.Range("I14").Value = 84
$ws.Range("J14").Value = 0.1064

$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 51
$ws.Range("M14").Value = 0.0893

$ws.Range("N14").Value = 4
$ws.Range("O14").Value = 46
$ws.Range("P14").Value = 0.08

$ws.Range("Q14").Value = 8
$ws.Range("R14").Value = 108
$ws.Range("S14").Value = 0.069

$ws.Range("T14").Value = 2
$ws.Range("U14").Value = 154
$ws.Range("V14").Value = 0.0128

# --- Total row (row 15): arrest_made, no_arrest_made, arrest_rate per year ---
$ws.Range("B15").Value = 37
$ws.Range("C15").Value = 288
$ws.Range("D15").Value = 0.1138

$ws.Range("E15").Value = 66
$ws.Range("F15").Value = 574
$ws.Range("G15").Value = 0.1031

$ws.Range("H15").Value = 73
$ws.Range("I15").Value = 842
$ws.Range("J15").Value = 0.0798

$ws.Range("K15").Value = 79
$ws.Range("L15").Value = 659
$ws.Range("M15").Value = 0.107

$ws.Range("N15").Value = 58
$ws.Range("O15").Value = 526
$ws.Range("P15").Value = 0.0993

$ws.Range("Q15").Value = 72
$ws.Range("R15").Value = 1308
$ws.Range("S15").Value = 0.0522

$ws.Range("T15").Value = 102
$ws.Range("U15").Value = 1697
$ws.Range("V15").Value = 0.0567
